$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lom3259 = "LOM3259 -  Materiais e Dispositivos Eletrônicos  (Indicação de Conjunto)`n"
$lom3234 = "LOM3234 -  Óptica Física  (Requisito)`n"

# The two requisite rows (24 and 25) swap their text: LOM3234 now comes
# first (row 24), LOM3259 moves to the second slot (row 25).
$ws.Range("B24").Value = $lom3234
$ws.Range("C24").Value = $lom3234
$ws.Range("B25").Value = $lom3259
$ws.Range("C25").Value = $lom3259
